$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows before row 347 (shifts old rows 347-395 down to 352-400)
$ws.Range("A347:A351").EntireRow.Insert()

# Row 347
$ws.Range("A347").Value = 6
$ws.Range("B347").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C347").Value = 'Metropolitana'
$ws.Range("D347").Value = 44505
$ws.Range("E347").Value = 13
$ws.Range("F347").Value = 100112003
$ws.Range("G347").Value = 'Ajo'
$ws.Range("H347").Value = 'Chino'
$ws.Range("I347").Value = '1a nueva(o)'
$ws.Range("J347").Value = 35000
$ws.Range("K347").Value = 1800
$ws.Range("L347").Value = 1800
$ws.Range("M347").Value = 1800
$ws.Range("N347").Value = '$/paquete 20 unidades (volumen en unidades)'
$ws.Range("O347").Value = 'Llay Llay'
$ws.Range("P347").Value = 90
$ws.Range("Q347").Value = 20
$ws.Range("R347").Value = 'Hortaliza'

# Row 348
$ws.Range("A348").Value = 6
$ws.Range("B348").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C348").Value = 'Metropolitana'
$ws.Range("D348").Value = 44505
$ws.Range("E348").Value = 13
$ws.Range("F348").Value = 100112003
$ws.Range("G348").Value = 'Ajo'
$ws.Range("H348").Value = 'Chino'
$ws.Range("I348").Value = '2a nueva(o)'
$ws.Range("J348").Value = 27000
$ws.Range("K348").Value = 1200
$ws.Range("L348").Value = 1200
$ws.Range("M348").Value = 1200
$ws.Range("N348").Value = '$/paquete 20 unidades (volumen en unidades)'
$ws.Range("O348").Value = 'Llay Llay'
$ws.Range("P348").Value = 60
$ws.Range("Q348").Value = 20
$ws.Range("R348").Value = 'Hortaliza'

# Row 349
$ws.Range("A349").Value = 6
$ws.Range("B349").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C349").Value = 'Metropolitana'
$ws.Range("D349").Value = 44505
$ws.Range("E349").Value = 13
$ws.Range("F349").Value = 100112003
$ws.Range("G349").Value = 'Ajo'
$ws.Range("H349").Value = 'Chino'
$ws.Range("I349").Value = '3a nueva (o)'
$ws.Range("J349").Value = 15000
$ws.Range("K349").Value = 800
$ws.Range("L349").Value = 800
$ws.Range("M349").Value = 800
$ws.Range("N349").Value = '$/paquete 20 unidades (volumen en unidades)'
$ws.Range("O349").Value = 'Llay Llay'
$ws.Range("P349").Value = 40
$ws.Range("Q349").Value = 20
$ws.Range("R349").Value = 'Hortaliza'

# Row 350
$ws.Range("A350").Value = 6
$ws.Range("B350").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C350").Value = 'Metropolitana'
$ws.Range("D350").Value = 44505
$ws.Range("E350").Value = 13
$ws.Range("F350").Value = 100112003
$ws.Range("G350").Value = 'Ajo'
$ws.Range("H350").Value = 'Chino'
$ws.Range("I350").Value = 'Extra nueva (o)'
$ws.Range("J350").Value = 33000
$ws.Range("K350").Value = 2400
$ws.Range("L350").Value = 2400
$ws.Range("M350").Value = 2400
$ws.Range("N350").Value = '$/paquete 20 unidades (volumen en unidades)'
$ws.Range("O350").Value = 'Llay Llay'
$ws.Range("P350").Value = 120
$ws.Range("Q350").Value = 20
$ws.Range("R350").Value = 'Hortaliza'

# Row 351
$ws.Range("A351").Value = 6
$ws.Range("B351").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C351").Value = 'Metropolitana'
$ws.Range("D351").Value = 44505
$ws.Range("E351").Value = 13
$ws.Range("F351").Value = 100112003
$ws.Range("G351").Value = 'Ajo'
$ws.Range("H351").Value = 'Chino'
$ws.Range("I351").Value = 'Primera'
$ws.Range("J351").Value = 2200
$ws.Range("K351").Value = 16500
$ws.Range("L351").Value = 17000
$ws.Range("M351").Value = 16773
$ws.Range("N351").Value = '$/caja 10 kilos'
$ws.Range("O351").Value = 'China'
$ws.Range("P351").Value = 1677
$ws.Range("Q351").Value = 10
$ws.Range("R351").Value = 'Hortaliza'

